$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("genotype")

# Fix the "cep290_unkown" typo -> "cep290_unknown" across the well-result grid.
# (Row 10's "failed" label cell, C10, is untouched and keeps its own text.)
$ws.Range("B2:M9").Value = "cep290_unknown"

# Bring the genotype sheet to the front and leave B2 selected, matching the
# author's final navigation state (this also clears tabSelected on whichever
# sheet -- "temperature" -- was active before).
$ws.Activate()
$ws.Range("B2").Select()
